$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.697.48"
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").Value = "'2.083.60"
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("D5").Value = "'228.23"
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("D7").Value = "'60.48"
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +1.92%  '
$ws.Range("D10").Value = "'0.0836"
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").Value = "'2.392.74"
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("E13").Value = '  +3.44%  '
$ws.Range("D14").Value = "'21.87"
$ws.Range("E14").Value = '  +3.59%  '
$ws.Range("E15").Value = '  +3.58%  '
$ws.Range("D16").Value = "'5.47"
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = "'2.086.16"
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").Value = "'38.648.50"
$ws.Range("E18").Value = '  +2.36%  '
$ws.Range("D19").Value = "'71.55"
$ws.Range("E19").Value = '  +3.01%  '
$ws.Range("D20").Value = "'6.01"
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("D22").Value = "'226.50"
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").Value = "'2.34"
$ws.Range("E25").Value = '  +2.29%  '
$ws.Range("D26").Value = "'170.81"
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("E28").Value = '  +7.10%  '
$ws.Range("E29").Value = '  +11.95%  '
$ws.Range("D30").Value = "'19.14"
$ws.Range("E30").Value = '  +1.95%  '
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("D32").Value = "'2.34"
$ws.Range("E32").Value = '  +4.22%  '
$ws.Range("E33").Value = '  +2.64%  '
$ws.Range("E34").Value = '  +4.54%  '
$ws.Range("D35").Value = "'0.0614"
$ws.Range("E35").Value = '  +2.04%  '
$ws.Range("E36").Value = '  +1.95%  '
$ws.Range("D37").Value = "'6.40"
$ws.Range("E37").Value = '  -2.50%  '
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = '  +2.17%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").Value = "'18.18"
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("D41").Value = "'1.539.55"
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("D42").Value = "'100.68"
$ws.Range("E42").Value = '  +2.91%  '
$ws.Range("D43").Value = "'0.0224"
$ws.Range("E43").Value = '  +3.84%  '
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("E45").Value = '  +2.02%  '
$ws.Range("D46").Value = "'7.68"
$ws.Range("E46").Value = '  +8.21%  '
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D48").Value = "'4.08"
$ws.Range("E48").Value = '  -2.48%  '
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").Value = "'2.281.40"
$ws.Range("E51").Value = '  +2.09%  '
